# Update "想去人数" (want-to-go count) figures for both the "展览" and
# "全部类型" sheets, which contain duplicate data tables.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1494
    $ws.Range("F7").Value = 118
    $ws.Range("F8").Value = 47
    $ws.Range("F9").Value = 281
}

$wb.Save()
